$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "description" column (C) values for the net-good-stats, swagger,
# damage+lower and damage+raise categories had their apostrophes doubled
# (e.g. "target's" -> "target''s"). Update the description text in place;
# all other data is unchanged.
$ws.Range("C4").Value = "No damage; lowers target''s stats or raises user''s stats"
$ws.Range("C7").Value = "No damage; inflicts status ailment; raises target''s stats"
$ws.Range("C8").Value = "Inflicts damage; lowers target''s stats"
$ws.Range("C9").Value = "Inflicts damage; raises user''s stats"
